$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10138.7145515933
$ws.Range("F2").Value = -3.39707875505625

$ws.Range("C3").Value = 10166.8739381789
$ws.Range("F3").Value = 243.4673752057

$ws.Range("C9").Value = 9714.35324426396
$ws.Range("F9").Value = 261.412673941831

$ws.Range("C10").Value = 9419.02060539191
$ws.Range("F10").Value = 249.107147322163

$ws.Range("C11").Value = 9419.77442110185
$ws.Range("F11").Value = 249.138556310076

$ws.Range("C12").Value = 8650.95608578513
$ws.Range("F12").Value = 217.104459005213

$ws.Range("C13").Value = 8408.88531461681
$ws.Range("F13").Value = 190.235187144572

$ws.Range("C14").Value = 8086.86611537369
$ws.Range("F14").Value = 176.305478656315

$ws.Range("C15").Value = 8955.7968205768
$ws.Range("F15").Value = 241.197099526731
